$d = $word.ActiveDocument

# Set ligatures on the Normal style's font (drives styles.xml Normal w:rPr)
$normal = $d.Styles(1)  # wdStyleNormal = -1 normally, but try by name too
